# Correction : if there is more places than babies, ensure that a baby has only one room
#
# The "rooms" sheet had the special "out" bookkeeping row immediately after
# the last real room (r10). Three brand-new "rea" rooms (r11, r12, r13) are
# inserted above that bookkeeping row so that the hospital has more rooms
# than babies - each baby can then be guaranteed a single room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rooms")

# Insert three blank rows at row 12 (just above the "out" row), pushing the
# "out" bookkeeping row (and its special shaded style) down to row 15.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# Populate the three newly inserted rows with the new rooms: they are brand
# new rooms (not part of the old room list), all allocated to the "rea"
# service, each with a capacity of 1.
$newRooms = @("r11", "r12", "r13")
for ($i = 0; $i -lt $newRooms.Length; $i++) {
    $r = 12 + $i
    $ws.Cells.Item($r, 1).Value = $newRooms[$i]   # all_rooms
    $ws.Cells.Item($r, 2).Value = "yes"            # new_rooms
    $ws.Cells.Item($r, 5).Value = "rea"            # new_rooms_service
    $ws.Cells.Item($r, 6).Value = "rea"            # old_rooms_service
    $ws.Cells.Item($r, 7).Value = 1                # rooms_capacities
}

# Make "rooms" the active sheet/tab and select the newly added block.
$ws.Activate()
$ws.Range("E12:G14").Select()
